$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

$ws.Range('D2').Value = '30.039.16'
$ws.Range('E2').Value = '  -0.21%  '
$ws.Range('D3').Value = '1.873.78'
$ws.Range('E3').Value = '  -2.53%  '
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('E5').Value = '  -3.62%  '
$ws.Range('D6').Value = '0.9989'
$ws.Range('E6').Value = '  -0.10%  '
$ws.Range('D7').Value = '0.5093'
$ws.Range('E7').Value = '  -3.03%  '
$ws.Range('D8').Value = '0.3953'
$ws.Range('E8').Value = '  -2.97%  '
$ws.Range('D9').Value = '0.08210'
$ws.Range('E9').Value = '  -3.94%  '
$ws.Range('D10').Value = '42.22'
$ws.Range('E10').Value = '  -2.84%  '
$ws.Range('E11').Value = '  -3.08%  '
$ws.Range('D12').Value = '23.95'
$ws.Range('E12').Value = '  +6.63%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = '6.320'
$ws.Range('E13').Value = '  -2.00%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.865.96'
$ws.Range('E14').Value = '  -2.84%  '
$ws.Range('D15').Value = '7.208'
$ws.Range('E15').Value = '  -2.87%  '
$ws.Range('D16').Value = '0.9999'
$ws.Range('E16').Value = '  -0.09%  '
$ws.Range('D17').Value = '92.23'
$ws.Range('E17').Value = '  -4.73%  '
$ws.Range('D18').Value = '0.00001087'
$ws.Range('E18').Value = '  -2.77%  '
$ws.Range('D19').Value = '0.06393'
$ws.Range('E19').Value = '  -4.82%  '
$ws.Range('D20').Value = '18.08'
$ws.Range('E20').Value = '  -1.50%  '
$ws.Range('D21').Value = '0.9991'
$ws.Range('E21').Value = '  -0.03%  '
$ws.Range('D22').Value = '30.038.84'
$ws.Range('E22').Value = '  -0.25%  '
$ws.Range('D23').Value = '5.851'
$ws.Range('E23').Value = '  -3.61%  '
$ws.Range('D24').Value = '11.14'
$ws.Range('E24').Value = '  -1.61%  '
$ws.Range('D25').Value = '2.173'
$ws.Range('E25').Value = '  -2.52%  '
$ws.Range('D26').Value = '2.083.30'
$ws.Range('E26').Value = '  -2.78%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = '21.17'
$ws.Range('E27').Value = '  -0.11%  '
$ws.Range('B28').Value = 'Monero'
$ws.Range('C28').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D28').Value = '160.64'
$ws.Range('E28').Value = '  +0.12%  '
$ws.Range('E29').Value = '  -9.42%  '
$ws.Range('E30').Value = '  -1.44%  '
$ws.Range('D31').Value = '1.076'
$ws.Range('E31').Value = '  -0.81%  '
$ws.Range('D32').Value = '0.1037'
$ws.Range('E32').Value = '  -2.18%  '
$ws.Range('D33').Value = '5.972'
$ws.Range('E33').Value = '  -2.67%  '
$ws.Range('D34').Value = '3.712'
$ws.Range('E34').Value = '  +1.83%  '
$ws.Range('D35').Value = '0.02443'
$ws.Range('E35').Value = '  -3.27%  '
$ws.Range('D36').Value = '5.261'
$ws.Range('E36').Value = '  +0.65%  '
$ws.Range('D37').Value = '0.06406'
$ws.Range('E37').Value = '  -3.23%  '
$ws.Range('D38').Value = '0.2151'
$ws.Range('E38').Value = '  -3.54%  '
$ws.Range('D39').Value = '1.181'
$ws.Range('E39').Value = '  -4.73%  '
$ws.Range('D40').Value = '8.558'
$ws.Range('E40').Value = '  -5.64%  '
$ws.Range('D41').Value = '11.44'
$ws.Range('E41').Value = '  -2.41%  '
$ws.Range('D42').Value = '0.6326'
$ws.Range('E42').Value = '  -3.84%  '
$ws.Range('D43').Value = '1.206'
$ws.Range('E43').Value = '  -3.08%  '
$ws.Range('D44').Value = '0.9983'
$ws.Range('E44').Value = '  -0.04%  '
$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D45').Value = '0.5926'
$ws.Range('E45').Value = '  -4.60%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').Value = '12.95'
$ws.Range('E46').Value = '  -3.26%  '
$ws.Range('D47').Value = '3.640'
$ws.Range('E47').Value = '  -4.01%  '
$ws.Range('D48').Value = '2.026'
$ws.Range('E48').Value = '  -3.33%  '
$ws.Range('D49').Value = '122.86'
$ws.Range('E49').Value = '  -2.01%  '
$ws.Range('D50').Value = '1.209'
$ws.Range('E50').Value = '  -3.23%  '
$ws.Range('D51').Value = '1.124'
$ws.Range('E51').Value = '  -2.96%  '
